$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Paragraph-internal run merges: replacing a span of text (that
#    happens to cross run boundaries) with the exact same text causes
#    Word to coalesce the runs into a single run on save - this mirrors
#    the "de-duplicated run" shape seen throughout the diff.
# ---------------------------------------------------------------------

$d.Content.Find.Execute(
    "Provide a faster tracing alternative than IDA Debugging, which can be very slow.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Provide a faster tracing alternative than IDA Debugging, which can be very slow.", 2)

$d.Content.Find.Execute(
    "CBASS execution PIN tracer (exetrace.dll) can track input at system call level so it capture much more comprehensive inputs from file and network than Win32 API level.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "CBASS execution PIN tracer (exetrace.dll) can track input at system call level so it capture much more comprehensive inputs from file and network than Win32 API level.", 2)

$d.Content.Find.Execute(
    "CBASS execution PIN tracer (exetrace.dll) uses binary encoding and instruction map to reduce redundancy in trace file, so the trace file is significantly smaller.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "CBASS execution PIN tracer (exetrace.dll) uses binary encoding and instruction map to reduce redundancy in trace file, so the trace file is significantly smaller.", 2)

$d.Content.Find.Execute(
    "PIN only supports x86 on Windows/Linux, and CBASS execution PIN tracer",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "PIN only supports x86 on Windows/Linux, and CBASS execution PIN tracer", 2)

$d.Content.Find.Execute(
    "CBASS PIN Tracer currently supports one file filter and no network filter.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "CBASS PIN Tracer currently supports one file filter and no network filter.", 2)

$d.Content.Find.Execute(
    "To avoid sending trace files back and forth between TREE GUI and target machine, a shared folder named TREE-TRACE is mapped to Z:\TREE-TRACE inside target VM. CBASS ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "To avoid sending trace files back and forth between TREE GUI and target machine, a shared folder named TREE-TRACE is mapped to Z:\TREE-TRACE inside target VM. CBASS ", 2)

# ---------------------------------------------------------------------
# 2) Append the new "PinAgent error" sentence after the Trace-Ready
#    sentence, at the end of that bullet paragraph.
# ---------------------------------------------------------------------

$rngTraceReady = $d.Content
$rngTraceReady.Find.Execute(
    " will send a Trace Ready message back to TREE GUI, so TREE GUI can notify user to move on to next stage.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngTraceReady.InsertAfter(" If an error is generated from PIN process, PinAgent will send an error message to GUI.")

# ---------------------------------------------------------------------
# 3) Sanity-test paragraph: remove the mid-word bookmark split and
#    merge "... traces f" / "or these cases are confirmed." back into
#    a single run of plain text (the _GoBack bookmark moves to the end
#    of the document - handled in step (6) below).
# ---------------------------------------------------------------------

$d.Content.Find.Execute(
    " on Windows XP SP3 and SP2, traces for these cases are confirmed.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " on Windows XP SP3 and SP2, traces for these cases are confirmed.", 2)

# ---------------------------------------------------------------------
# 4) TODO list: append "[DONE]" after the first item's "(Xing)" marker.
# ---------------------------------------------------------------------

$rngDone1 = $d.Content
$rngDone1.Find.Execute(
    "Update installation process to include PIN tracer components.(Xing)",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngDone1.InsertAfter("[DONE]")

# ---------------------------------------------------------------------
# 5) TODO list: append "[DONE}" after the Tracer-GUI item.
# ---------------------------------------------------------------------

$rngDone2 = $d.Content
$rngDone2.Find.Execute(
    "Update Tracer GUI for new configurations and enhance usability.(Xing)",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngDone2.InsertAfter("[DONE}")

# ---------------------------------------------------------------------
# 6) TODO list: append "[Needs Integration]" after the last item, then
#    re-create the _GoBack bookmark at the very end of the document.
#    A temporary marker character is used to land the (zero-width)
#    bookmark exactly at the end of the paragraph's text, then that
#    marker is deleted while the bookmark stays in place.
# ---------------------------------------------------------------------

$rngLast = $d.Content
$rngLast.Find.Execute(
    "Develop new PIN-trace parser and integrate with Taint Analyzer. (Nathan)",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngLast.InsertAfter("[Needs Integration]")
$rngLast.InsertAfter("Z")

$bmRange = $d.Range($rngLast.End - 1, $rngLast.End)
$d.Bookmarks.Add("_GoBack", $bmRange)

$markerRange = $d.Range($rngLast.End - 1, $rngLast.End)
$markerRange.Text = ""
